# Generate Report for Handback
#
# Refresh the handback-status report: the generator re-ran and produced
# newer "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 100bebab...md file on the per-language sheets, and a
# newer "Latest HO Xliff Generate Date" roll-up for the cd7bd2de...md file
# on the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to cd7bd2de-be0d-4ed1-8036-d90dd773d8e9.md
$overview.Range("G3").Value = "2016-09-01 18:57:20"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 2 corresponds to 100bebab-988a-4514-be6c-9bd0731fed9d.md
$zhcn.Range("H2").Value = "2016-09-01 18:57:10"
$zhcn.Range("K2").Value = "2016-09-01 18:57:40"

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# Row 2 corresponds to 100bebab-988a-4514-be6c-9bd0731fed9d.md
$dede.Range("H2").Value = "2016-09-01 18:57:20"
$dede.Range("K2").Value = "2016-09-01 18:57:48"
